$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '68.117.29'
$ws.Range('E2').Value = '  -0.97%  '
$ws.Range('D3').Value = '2.643.64'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = "'596.52"
$ws.Range('E5').Value = '  -0.76%  '
$ws.Range('D6').Value = "'156.70"
$ws.Range('E6').Value = '  +1.04%  '
$ws.Range('E8').Value = '  -0.33%  '
$ws.Range('E9').Value = '  +2.40%  '
$ws.Range('E10').Value = '  -1.21%  '
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').Value = "'28.05"
$ws.Range('E13').Value = '  +0.93%  '
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').Value = '3.125.34'
$ws.Range('E15').Value = '  +0.01%  '
$ws.Range('D16').Value = '68.165.72'
$ws.Range('E16').Value = '  -0.62%  '
$ws.Range('D17').Value = '2.637.83'
$ws.Range('E17').Value = '  -0.08%  '
$ws.Range('D18').Value = "'11.36"
$ws.Range('E18').Value = '  -0.78%  '
$ws.Range('D19').Value = "'362.41"
$ws.Range('E19').Value = '  -1.44%  '
$ws.Range('D20').Value = "'7.40"
$ws.Range('E20').Value = '  -0.93%  '
$ws.Range('D21').Value = "'4.43"
$ws.Range('E21').Value = '  +3.85%  '
$ws.Range('E22').Value = '  -1.53%  '
$ws.Range('E23').Value = '  -1.87%  '
$ws.Range('D24').Value = "'75.07"
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D26').Value = "'9.68"
$ws.Range('E26').Value = '  -3.22%  '
$ws.Range('E28').Value = '  -2.27%  '
$ws.Range('D29').Value = "'0.998"
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('D30').Value = "'554.90"
$ws.Range('E30').Value = '  -5.00%  '
$ws.Range('D31').Value = "'8.01"
$ws.Range('E31').Value = '  +0.10%  '
$ws.Range('E32').Value = '  -1.49%  '
$ws.Range('E33').Value = '  -0.55%  '
$ws.Range('E34').Value = '  -1.39%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('E36').Value = '  +0.49%  '
$ws.Range('D37').Value = "'161.06"
$ws.Range('E37').Value = '  +1.14%  '
$ws.Range('D38').Value = "'19.39"
$ws.Range('E38').Value = '  +0.32%  '
$ws.Range('E39').Value = '  +1.04%  '
$ws.Range('E40').Value = '  -2.77%  '
$ws.Range('D41').Value = "'5.33"
$ws.Range('E41').Value = '  -1.26%  '
$ws.Range('D42').Value = '0.0₆0339'
$ws.Range('E42').Value = '  +6.45%  '
$ws.Range('D43').Value = "'17.78"
$ws.Range('E43').Value = '  +0.44%  '
$ws.Range('D44').Value = "'2.61"
$ws.Range('E44').Value = '  -2.64%  '
$ws.Range('E45').Value = '  +0.03%  '
$ws.Range('D46').Value = "'40.40"
$ws.Range('E46').Value = '  -0.34%  '
$ws.Range('D47').Value = "'158.86"
$ws.Range('E47').Value = '  +1.69%  '
$ws.Range('E48').Value = '  -0.06%  '
$ws.Range('D49').Value = "'22.05"
$ws.Range('E49').Value = '  -0.26%  '
$ws.Range('E50').Value = '  -0.94%  '
$ws.Range('D51').Value = "'0.0784"
$ws.Range('E51').Value = '  +0.24%  '

Write-Host "Done updating cells."
